$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated financial figures (Income Statement, Balance Sheet, Cash Flow Statement)
$ws.Range("D8").Value = 258800
$ws.Range("E8").Value = 74900
$ws.Range("D9").Value = 107300
$ws.Range("E9").Value = 85800
$ws.Range("D10").Value = 151600
$ws.Range("E10").Value = -10800
$ws.Range("D12").Value = 19000
$ws.Range("E12").Value = 4300
$ws.Range("D15").Value = 300
$ws.Range("D17").Value = 347200
$ws.Range("E17").Value = 117400
$ws.Range("D18").Value = -88400
$ws.Range("E18").Value = -42500
$ws.Range("D20").Value = 10500
$ws.Range("E20").Value = -900
$ws.Range("D21").Value = -77600
$ws.Range("E21").Value = -43200
$ws.Range("D23").Value = -77900
$ws.Range("E23").Value = -43300
$ws.Range("D26").Value = -77900
$ws.Range("E26").Value = -43300
$ws.Range("D27").Value = -74000
$ws.Range("E27").Value = -47800
$ws.Range("D32").Value = -10500
$ws.Range("E32").Value = 900
$ws.Range("D33").Value = -74000
$ws.Range("E33").Value = -47800
$ws.Range("D35").Value = -74000
$ws.Range("E35").Value = -47800
$ws.Range("D41").Value = 453900
$ws.Range("E41").Value = 195900
$ws.Range("D42").Value = 7400
$ws.Range("E42").Value = 43000
$ws.Range("D43").Value = 88100
$ws.Range("E43").Value = 16700
$ws.Range("D45").Value = 1400400
$ws.Range("E45").Value = 4600
$ws.Range("D46").Value = 1949800
$ws.Range("E46").Value = 260200
$ws.Range("D47").Value = 24800
$ws.Range("D48").Value = 1400
$ws.Range("D54").Value = 1976000
$ws.Range("E54").Value = 262800
$ws.Range("D57").Value = 1460100
$ws.Range("E57").Value = 167500
$ws.Range("D59").Value = 337000
$ws.Range("E59").Value = 42400
$ws.Range("D60").Value = 1797200
$ws.Range("E60").Value = 209900
$ws.Range("D66").Value = 1797200
$ws.Range("E66").Value = 209900
$ws.Range("D70").Value = 326000
$ws.Range("E70").Value = 116200
$ws.Range("D72").Value = -152900
$ws.Range("E72").Value = -70100
$ws.Range("D76").Value = -147200
$ws.Range("E76").Value = -63300
$ws.Range("D81").Value = -74000
$ws.Range("E81").Value = -47800
$ws.Range("D89").Value = 46800
$ws.Range("E89").Value = 130600
$ws.Range("D94").Value = 10600
$ws.Range("E94").Value = -45600
$ws.Range("D100").Value = 207600
$ws.Range("E100").Value = 72200
$ws.Range("D101").Value = -7100
$ws.Range("D102").Value = 258000
$ws.Range("E102").Value = 160200
